$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Classrooms" header in column E (row 1)
$ws.Range("E1").Value = "Classrooms"

# Set the width of the new column E (close to 20.33203125 characters)
$ws.Columns.Item(5).ColumnWidth = 19.5

# Touch E2 briefly so the sheet's used range / row span extends to column E,
# then clear the value back out (keeps row 2 empty, matching the target row).
$ws.Range("E2").Value = "tmp"
$ws.Range("E2").ClearContents()

# Update the selected cell to match the new state
$ws.Range("C4").Select()
